$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update model labels in column A (row order permutation)
$ws.Range("A2").Value = "model_12_5_0"
$ws.Range("A3").Value = "model_12_5_22"
$ws.Range("A4").Value = "model_12_5_21"
$ws.Range("A5").Value = "model_12_5_20"
$ws.Range("A6").Value = "model_12_5_19"
$ws.Range("A7").Value = "model_12_5_18"
$ws.Range("A8").Value = "model_12_5_17"
$ws.Range("A9").Value = "model_12_5_16"
$ws.Range("A10").Value = "model_12_5_15"
$ws.Range("A11").Value = "model_12_5_14"
$ws.Range("A12").Value = "model_12_5_13"
$ws.Range("A13").Value = "model_12_5_23"
$ws.Range("A14").Value = "model_12_5_12"
$ws.Range("A15").Value = "model_12_5_10"
$ws.Range("A16").Value = "model_12_5_9"
$ws.Range("A17").Value = "model_12_5_8"
$ws.Range("A18").Value = "model_12_5_7"
$ws.Range("A19").Value = "model_12_5_6"
$ws.Range("A20").Value = "model_12_5_5"
$ws.Range("A21").Value = "model_12_5_4"
$ws.Range("A22").Value = "model_12_5_3"
$ws.Range("A23").Value = "model_12_5_2"
$ws.Range("A24").Value = "model_12_5_1"
$ws.Range("A25").Value = "model_12_5_11"
$ws.Range("A26").Value = "model_12_5_24"

# Update metric columns B..Q - every row now shares the same metric values
$ws.Range("B2:B26").Value = [double]"0.9994384710386303"
$ws.Range("C2:C26").Value = [double]"0.9988584094093185"
$ws.Range("D2:D26").Value = [double]"0.9999999999998973"
$ws.Range("E2:E26").Value = [double]"0.9999994636804617"
$ws.Range("F2:F26").Value = [double]"0.999999664704229"
$ws.Range("G2:G26").Value = [double]"0.0005241626371659963"
$ws.Range("H2:H26").Value = [double]"0.001065624706365756"
$ws.Range("I2:I26").Value = [double]"5.09579840977513e-14"
$ws.Range("J2:J26").Value = [double]"4.50933016801211e-07"
$ws.Range("K2:K26").Value = [double]"2.254665338795976e-07"
$ws.Range("L2:L26").Value = [double]"0.0009694867605113856"
$ws.Range("M2:M26").Value = [double]"0.02289459842770771"
$ws.Range("N2:N26").Value = [double]"1.002695339014575"
$ws.Range("O2:O26").Value = [double]"0.02386926888112749"
$ws.Range("P2:P26").Value = [double]"73.107417091058"
$ws.Range("Q2:Q26").Value = [double]"108.4548160122358"
